$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Celdas Naranjas" (column C) counts following the recount ("contar") pass.
$ws.Range("C5").Value  = 0
$ws.Range("C8").Value  = 0
$ws.Range("C16").Value = 0
$ws.Range("C17").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("C24").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("C30").Value = 0
$ws.Range("C31").Value = 0
$ws.Range("C32").Value = 0
$ws.Range("C33").Value = 0
$ws.Range("C34").Value = 0
$ws.Range("C35").Value = 2
$ws.Range("C36").Value = 0
$ws.Range("C37").Value = 0
$ws.Range("C38").Value = 0
$ws.Range("C39").Value = 0
$ws.Range("C40").Value = 0
$ws.Range("C41").Value = 0
$ws.Range("C42").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("C44").Value = 0
$ws.Range("C45").Value = 0
$ws.Range("C46").Value = 2
$ws.Range("C47").Value = 0
$ws.Range("C49").Value = 0
$ws.Range("C50").Value = 22
$ws.Range("C60").Value = 97
